$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column J = "saptamana 8" - mark 1 exercise completed for these students
$rows = @(6, 9, 12, 13, 18, 19, 21)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 10).Value = 1
}

# Update active selection to reflect the last edited location
$ws.Range("L13").Select()
